# Update "想去人数" (interest count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 220
$ws1.Range("F5").Value = 1439
$ws1.Range("F7").Value = 604
$ws1.Range("F8").Value = 124
$ws1.Range("F9").Value = 585
$ws1.Range("F10").Value = 31
$ws1.Range("F11").Value = 672

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 4
$ws2.Range("F11").Value = 9

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6287
$ws3.Range("F3").Value = 778

# Sheet "全部类型" (All Types) - aggregated view of the above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6287
$ws4.Range("F3").Value = 778
$ws4.Range("F11").Value = 220
$ws4.Range("F12").Value = 4
$ws4.Range("F15").Value = 1439
$ws4.Range("F19").Value = 9
$ws4.Range("F20").Value = 604
$ws4.Range("F22").Value = 124
$ws4.Range("F23").Value = 585
$ws4.Range("F24").Value = 31
$ws4.Range("F26").Value = 672
